$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2386831275720165
$ws.Range("C2").Value = 0.4938271604938271
$ws.Range("J2").Value = 0.0205761316872428
$ws.Range("P2").Value = 0.1646090534979424
$ws.Range("S2").Value = 0.08230452674897119
$ws.Range("B3").Value = 0.01680672268907563
$ws.Range("C3").Value = 0.02521008403361345
$ws.Range("J3").Value = 0.02521008403361345
$ws.Range("P3").Value = 0.7058823529411765
$ws.Range("S3").Value = 0.226890756302521
$ws.Range("J4").Value = 0.06896551724137931
$ws.Range("P4").Value = 0.6551724137931034
$ws.Range("S4").Value = 0.2758620689655172
$ws.Range("B6").Value = 0.08235294117647059
$ws.Range("F6").Value = 0.07647058823529412
$ws.Range("J6").Value = 0.2
$ws.Range("O6").Value = 0.02352941176470588
$ws.Range("Q6").Value = 0.1176470588235294
$ws.Range("R6").Value = 0.1
$ws.Range("S6").Value = 0.4
$ws.Range("B7").Value = 0.07777777777777778
$ws.Range("D7").Value = 0.005555555555555556
$ws.Range("F7").Value = 0.04444444444444445
$ws.Range("J7").Value = 0.1611111111111111
$ws.Range("Q7").Value = 0.1722222222222222
$ws.Range("R7").Value = 0.09444444444444444
$ws.Range("S7").Value = 0.4444444444444444
$ws.Range("B8").Value = 0.07399103139013453
$ws.Range("D8").Value = 0.02017937219730942
$ws.Range("F8").Value = 0.06278026905829596
$ws.Range("J8").Value = 0.08744394618834081
$ws.Range("O8").Value = 0.004484304932735426
$ws.Range("Q8").Value = 0.1883408071748879
$ws.Range("R8").Value = 0.08295964125560538
$ws.Range("S8").Value = 0.4798206278026906
$ws.Range("B9").Value = 0.07142857142857142
$ws.Range("D9").Value = 0.01948051948051948
$ws.Range("F9").Value = 0.04545454545454546
$ws.Range("J9").Value = 0.09740259740259741
$ws.Range("O9").Value = 0.01948051948051948
$ws.Range("Q9").Value = 0.2012987012987013
$ws.Range("R9").Value = 0.1363636363636364
$ws.Range("S9").Value = 0.4090909090909091
$ws.Range("B10").Value = 0.09819819819819819
$ws.Range("D10").Value = 0.01621621621621622
$ws.Range("E10").Value = 0.0009009009009009009
$ws.Range("F10").Value = 0.05765765765765766
$ws.Range("J10").Value = 0.1189189189189189
$ws.Range("O10").Value = 0.006306306306306306
$ws.Range("Q10").Value = 0.2135135135135135
$ws.Range("R10").Value = 0.07657657657657657
$ws.Range("S10").Value = 0.4117117117117117
$ws.Range("G11").Value = 0.1684210526315789
$ws.Range("J11").Value = 0.08421052631578947
$ws.Range("K11").Value = 0.2385964912280702
$ws.Range("L11").Value = 0.5017543859649123
$ws.Range("S11").Value = 0.007017543859649123
$ws.Range("G12").Value = 0.7762237762237763
$ws.Range("J12").Value = 0.1888111888111888
$ws.Range("K12").Value = 0.01398601398601399
$ws.Range("S12").Value = 0.02097902097902098
$ws.Range("G13").Value = 0.84375
$ws.Range("J13").Value = 0.125
$ws.Range("S13").Value = 0.03125
$ws.Range("F15").Value = 0.0108695652173913
$ws.Range("H15").Value = 0.2173913043478261
$ws.Range("I15").Value = 0.05978260869565218
$ws.Range("J15").Value = 0.4076086956521739
$ws.Range("K15").Value = 0.07065217391304347
$ws.Range("M15").Value = 0.005434782608695652
$ws.Range("O15").Value = 0.04347826086956522
$ws.Range("S15").Value = 0.1847826086956522
$ws.Range("F16").Value = 0.02857142857142857
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.07857142857142857
$ws.Range("J16").Value = 0.3928571428571428
$ws.Range("K16").Value = 0.06428571428571428
$ws.Range("M16").Value = 0.02142857142857143
$ws.Range("N16").Value = 0.007142857142857143
$ws.Range("O16").Value = 0.09285714285714286
$ws.Range("S16").Value = 0.1142857142857143
$ws.Range("F17").Value = 0.007462686567164179
$ws.Range("H17").Value = 0.1940298507462687
$ws.Range("I17").Value = 0.08208955223880597
$ws.Range("J17").Value = 0.4402985074626866
$ws.Range("K17").Value = 0.0945273631840796
$ws.Range("M17").Value = 0.007462686567164179
$ws.Range("O17").Value = 0.07960199004975124
$ws.Range("S17").Value = 0.0945273631840796
$ws.Range("F18").Value = 0.02808988764044944
$ws.Range("H18").Value = 0.1573033707865168
$ws.Range("I18").Value = 0.08426966292134831
$ws.Range("J18").Value = 0.4044943820224719
$ws.Range("K18").Value = 0.0898876404494382
$ws.Range("M18").Value = 0.02808988764044944
$ws.Range("O18").Value = 0.06741573033707865
$ws.Range("S18").Value = 0.1404494382022472
$ws.Range("F19").Value = 0.01279863481228669
$ws.Range("H19").Value = 0.2363481228668942
$ws.Range("I19").Value = 0.07337883959044368
$ws.Range("J19").Value = 0.3686006825938566
$ws.Range("K19").Value = 0.113481228668942
$ws.Range("M19").Value = 0.01706484641638225
$ws.Range("O19").Value = 0.06655290102389079
$ws.Range("S19").Value = 0.1117747440273038
